$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header and update row 2 values
$ws.Range("D1").Value = "Another animal"

$ws.Range("A2").Value = 41
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 3
